$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 3237.7827
$ws.Range("J112").Value = 3298.5908
$ws.Range("L112").Value = 9895.7724
$ws.Range("N112").Value = -12111.7724
$ws.Range("H129").Value = 1176.921
$ws.Range("I129").Value = 711.8889
$ws.Range("J129").Value = 1321.2413
$ws.Range("K129").Value = 2135.6667
$ws.Range("L129").Value = 3963.7239
$ws.Range("M129").Value = 2864.3333
$ws.Range("N129").Value = -13963.7239
$ws.Range("H137").Value = 1232.9584
$ws.Range("I137").Value = 1075.7646
$ws.Range("J137").Value = 1614.7142
$ws.Range("K137").Value = 3227.2938
$ws.Range("L137").Value = 4844.142599999999
$ws.Range("M137").Value = -677.2937999999999
$ws.Range("N137").Value = -9944.142599999999
$ws.Range("H138").Value = 2497.6
$ws.Range("I138").Value = 2670
$ws.Range("J138").Value = 2440.1333
$ws.Range("K138").Value = 8010
$ws.Range("L138").Value = 7320.3999
$ws.Range("M138").Value = -2870
$ws.Range("N138").Value = -17600.3999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 27999.5
$ws.Range("J55").Value = 27999.5
$ws.Range("L55").Value = 27999.5
$ws.Range("N55").Value = -28629.5
$ws.Range("H74").Value = 1277.8125
$ws.Range("I74").Value = 1210.25
$ws.Range("J74").Value = 1300.3334
$ws.Range("K74").Value = 1210.25
$ws.Range("L74").Value = 1300.3334
$ws.Range("M74").Value = -336.25
$ws.Range("N74").Value = -3048.3334
$ws.Range("H77").Value = 1277.8125
$ws.Range("I77").Value = 1210.25
$ws.Range("J77").Value = 1300.3334
$ws.Range("K77").Value = 6051.25
$ws.Range("L77").Value = 6501.666999999999
$ws.Range("M77").Value = -1683.25
$ws.Range("N77").Value = -15237.667
$ws.Range("H80").Value = 21665.666
$ws.Range("J80").Value = 21665.666
$ws.Range("L80").Value = 21665.666
$ws.Range("N80").Value = -23661.666
$ws.Range("H83").Value = 21665.666
$ws.Range("J83").Value = 21665.666
$ws.Range("L83").Value = 64996.99800000001
$ws.Range("N83").Value = -74980.99800000001
$ws.Range("H113").Value = 35756
$ws.Range("J113").Value = 35756
$ws.Range("L113").Value = 35756
$ws.Range("N113").Value = -44434

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2869.2222
$ws.Range("I134").Value = 2485.8462
$ws.Range("J134").Value = 3866
$ws.Range("K134").Value = 7457.5386
$ws.Range("L134").Value = 11598
$ws.Range("M134").Value = -4922.5386
$ws.Range("N134").Value = -16668

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5439.8687
$ws.Range("I31").Value = 1240.1818
$ws.Range("J31").Value = 11214.4375
$ws.Range("K31").Value = 1240.1818
$ws.Range("L31").Value = 11214.4375
$ws.Range("M31").Value = -945.1818000000001
$ws.Range("N31").Value = -11804.4375
$ws.Range("H34").Value = 5439.8687
$ws.Range("I34").Value = 1240.1818
$ws.Range("J34").Value = 11214.4375
$ws.Range("K34").Value = 1240.1818
$ws.Range("L34").Value = 11214.4375
$ws.Range("M34").Value = -1038.1818
$ws.Range("N34").Value = -11618.4375
$ws.Range("H132").Value = 20837204
$ws.Range("I132").Value = 3792.4
$ws.Range("J132").Value = 55559556
$ws.Range("K132").Value = 11377.2
$ws.Range("L132").Value = 166678668
$ws.Range("M132").Value = -8847.200000000001
$ws.Range("N132").Value = -166683728

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 599.18866
$ws.Range("I5").Value = 475.6154
$ws.Range("J5").Value = 943.4286
$ws.Range("K5").Value = 1426.8462
$ws.Range("L5").Value = 2830.2858
$ws.Range("M5").Value = -1314.8462
$ws.Range("N5").Value = -3054.2858
$ws.Range("H68").Value = 1437.2697
$ws.Range("J68").Value = 1596.2428
$ws.Range("L68").Value = 4788.7284
$ws.Range("N68").Value = -6410.7284
$ws.Range("H71").Value = 1437.2697
$ws.Range("J71").Value = 1596.2428
$ws.Range("L71").Value = 14366.1852
$ws.Range("N71").Value = -22478.1852
$ws.Range("H107").Value = 1389.64
$ws.Range("I107").Value = 268.92105
$ws.Range("J107").Value = 2540.6487
$ws.Range("K107").Value = 806.76315
$ws.Range("L107").Value = 7621.946100000001
$ws.Range("M107").Value = 1113.23685
$ws.Range("N107").Value = -11461.9461
$ws.Range("H113").Value = 880.6829
$ws.Range("I113").Value = 513.52
$ws.Range("J113").Value = 1454.375
$ws.Range("K113").Value = 1540.56
$ws.Range("L113").Value = 4363.125
$ws.Range("M113").Value = 629.4400000000001
$ws.Range("N113").Value = -8703.125
$ws.Range("H135").Value = 599.18866
$ws.Range("I135").Value = 475.6154
$ws.Range("J135").Value = 943.4286
$ws.Range("K135").Value = 4280.5386
$ws.Range("L135").Value = 8490.857399999999
$ws.Range("M135").Value = -1745.5386
$ws.Range("N135").Value = -13560.8574

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 49145
$ws.Range("J42").Value = 49145
$ws.Range("L42").Value = 49145
$ws.Range("N42").Value = -50115
$ws.Range("H102").Value = 696
$ws.Range("I102").Value = 696
$ws.Range("K102").Value = 696
$ws.Range("M102").Value = 926
$ws.Range("H115").Value = 49145
$ws.Range("J115").Value = 49145
$ws.Range("L115").Value = 49145
$ws.Range("N115").Value = -51495
$ws.Range("H126").Value = 3002
$ws.Range("I126").Value = 3002
$ws.Range("K126").Value = 9006
$ws.Range("M126").Value = -6536
$ws.Range("H132").Value = 3931.7144
$ws.Range("I132").Value = 3312
$ws.Range("J132").Value = 4179.6
$ws.Range("K132").Value = 9936
$ws.Range("L132").Value = 12538.8
$ws.Range("M132").Value = -7406
$ws.Range("N132").Value = -17598.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 6411561
$ws.Range("I136").Value = 1260.4445
$ws.Range("J136").Value = 20834738
$ws.Range("K136").Value = 3781.3335
$ws.Range("L136").Value = 62504214
$ws.Range("M136").Value = -1231.3335
$ws.Range("N136").Value = -62509314

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1625.4584
$ws.Range("I122").Value = 1375.375
$ws.Range("J122").Value = 2125.625
$ws.Range("K122").Value = 4126.125
$ws.Range("L122").Value = 6376.875
$ws.Range("M122").Value = -1676.125
$ws.Range("N122").Value = -11276.875
$ws.Range("H132").Value = 22439876
$ws.Range("I132").Value = 6251
$ws.Range("K132").Value = 18753
$ws.Range("M132").Value = -16223
$ws.Range("H136").Value = 3359.7188
$ws.Range("I136").Value = 3418.5881
$ws.Range("J136").Value = 3293
$ws.Range("K136").Value = 10255.7643
$ws.Range("L136").Value = 9879
$ws.Range("M136").Value = -7705.764299999999
$ws.Range("N136").Value = -14979

Write-Output "Applied 167 cell updates across 8 sheets"